$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "CHLORAMPHENICOL-2.0- UNSPECIFIED"
$ws.Range("B13").Value = "RIFAMPICIN-2.0 - UNSPECIFIED"
$ws.Range("B19").Value = "SPECTINOMYCIN-6.0 - UNSPECIFIED"

$ws.Range("B20").Select()
